{"js": "// Map of old cell text -> new cell text, as described by the diff.\nconst replacements = {\n  \"496\u00f74=124, 0\": \"681\u00f74=170, 1\",\n  \"762\u00f78=95, 2\": \"681\u00f77=97, 2\",\n  \"176\u00f72=88, 0\": \"463\u00f72=231, 1\",\n  \"210\u00f79=23, 3\": \"504\u00f72=252, 0\",\n  \"338\u00f72=169, 0\": \"525\u00f79=58, 3\",\n  \"564\u00f78=70, 4\": \"404\u00f77=57, 5\",\n  \"665\u00f79=73, 8\": \"649\u00f77=92, 5\",\n  \"746\u00f79=82, 8\": \"408\u00f73=136, 0\",\n  \"820\u00f74=205, 0\": \"916\u00f75=183, 1\",\n  \"101\u00f74=25, 1\": \"785\u00f74=196, 1\",\n  \"263\u00f77=37, 4\": \"689\u00f72=344, 1\",\n  \"546\u00f72=273, 0\": \"577\u00f74=144, 1\",\n  \"821\u00f76=136, 5\": \"403\u00f73=134, 1\",\n  \"688\u00f79=76, 4\": \"850\u00f77=121, 3\",\n  \"267\u00f75=53, 2\": \"937\u00f76=156, 1\",\n  \"356\u00f74=89, 0\": \"876\u00f78=109, 4\",\n  \"932\u00f79=103, 5\": \"988\u00f75=197, 3\",\n  \"491\u00f76=81, 5\": \"620\u00f79=68, 8\",\n  \"859\u00f73=286, 1\": \"565\u00f74=141, 1\",\n  \"267\u00f73=89, 0\": \"700\u00f78=87, 4\",\n  \"284\u00f76=47, 2\": \"357\u00f77=51, 0\",\n  \"374\u00f73=124, 2\": \"864\u00f75=172, 4\",\n  \"251\u00f79=27, 8\": \"235\u00f75=47, 0\",\n  \"742\u00f73=247, 1\": \"490\u00f73=163, 1\",\n  \"332\u00f74=83, 0\": \"744\u00f79=82, 6\",\n};\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const current = para.text;\n  if (Object.prototype.hasOwnProperty.call(replacements, current)) {\n    para.insertText(replacements[current], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Map of old cell text -> new cell text, as described by the diff.\n$replacements = [ordered]@{\n  \"496\u00f74=124, 0\" = \"681\u00f74=170, 1\"\n  \"762\u00f78=95, 2\"  = \"681\u00f77=97, 2\"\n  \"176\u00f72=88, 0\"  = \"463\u00f72=231, 1\"\n  \"210\u00f79=23, 3\"  = \"504\u00f72=252, 0\"\n  \"338\u00f72=169, 0\" = \"525\u00f79=58, 3\"\n  \"564\u00f78=70, 4\"  = \"404\u00f77=57, 5\"\n  \"665\u00f79=73, 8\"  = \"649\u00f77=92, 5\"\n  \"746\u00f79=82, 8\"  = \"408\u00f73=136, 0\"\n  \"820\u00f74=205, 0\" = \"916\u00f75=183, 1\"\n  \"101\u00f74=25, 1\"  = \"785\u00f74=196, 1\"\n  \"263\u00f77=37, 4\"  = \"689\u00f72=344, 1\"\n  \"546\u00f72=273, 0\" = \"577\u00f74=144, 1\"\n  \"821\u00f76=136, 5\" = \"403\u00f73=134, 1\"\n  \"688\u00f79=76, 4\"  = \"850\u00f77=121, 3\"\n  \"267\u00f75=53, 2\"  = \"937\u00f76=156, 1\"\n  \"356\u00f74=89, 0\"  = \"876\u00f78=109, 4\"\n  \"932\u00f79=103, 5\" = \"988\u00f75=197, 3\"\n  \"491\u00f76=81, 5\"  = \"620\u00f79=68, 8\"\n  \"859\u00f73=286, 1\" = \"565\u00f74=141, 1\"\n  \"267\u00f73=89, 0\"  = \"700\u00f78=87, 4\"\n  \"284\u00f76=47, 2\"  = \"357\u00f77=51, 0\"\n  \"374\u00f73=124, 2\" = \"864\u00f75=172, 4\"\n  \"251\u00f79=27, 8\"  = \"235\u00f75=47, 0\"\n  \"742\u00f73=247, 1\" = \"490\u00f73=163, 1\"\n  \"332\u00f74=83, 0\"  = \"744\u00f79=82, 6\"\n}\n\nforeach ($oldText in $replacements.Keys) {\n  $newText = $replacements[$oldText]\n  $range = $d.Content\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
